# Refresh cached Asura-server profit figures (market board snapshot) across
# the per-job leve-profit tables. Each block below updates one row's
# currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) cells.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 17
$ws_ALC.Range("H17").Value = 77937.69500000001
$ws_ALC.Range("J17").Value = 77937.69500000001
$ws_ALC.Range("L17").Value = 233813.085
$ws_ALC.Range("N17").Value = -234149.085

# Row 28
$ws_ALC.Range("H28").Value = 1387.6666
$ws_ALC.Range("I28").Value = 2243.4
$ws_ALC.Range("J28").Value = 318
$ws_ALC.Range("K28").Value = 2243.4
$ws_ALC.Range("L28").Value = 318
$ws_ALC.Range("M28").Value = -1758.4
$ws_ALC.Range("N28").Value = -1288

# Row 107
$ws_ALC.Range("H107").Value = 707.125
$ws_ALC.Range("I107").Value = 400.16666
$ws_ALC.Range("J107").Value = 891.3
$ws_ALC.Range("K107").Value = 400.16666
$ws_ALC.Range("L107").Value = 891.3
$ws_ALC.Range("M107").Value = 1519.83334
$ws_ALC.Range("N107").Value = -4731.3

# Row 111
$ws_ALC.Range("H111").Value = 2058.182
$ws_ALC.Range("I111").Value = 1788.2
$ws_ALC.Range("K111").Value = 5364.6
$ws_ALC.Range("M111").Value = -2297.6

# Row 113
$ws_ALC.Range("H113").Value = 2149.5833
$ws_ALC.Range("I113").Value = 1571.7646
$ws_ALC.Range("J113").Value = 3552.8572
$ws_ALC.Range("K113").Value = 1571.7646
$ws_ALC.Range("L113").Value = 3552.8572
$ws_ALC.Range("M113").Value = 1682.2354
$ws_ALC.Range("N113").Value = -10060.8572

# --- ARM sheet ---
$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 2
$ws_ARM.Range("H2").Value = 38656.63
$ws_ARM.Range("I2").Value = 1138.7
$ws_ARM.Range("J2").Value = 60726
$ws_ARM.Range("K2").Value = 1138.7
$ws_ARM.Range("L2").Value = 60726
$ws_ARM.Range("M2").Value = -1025.7
$ws_ARM.Range("N2").Value = -60952

# Row 7
$ws_ARM.Range("H7").Value = 50000
$ws_ARM.Range("J7").Value = 50000
$ws_ARM.Range("L7").Value = 50000
$ws_ARM.Range("N7").Value = -50228

# Row 38
$ws_ARM.Range("H38").Value = 2906.3333
$ws_ARM.Range("I38").Value = 2906.3333
$ws_ARM.Range("K38").Value = 2906.3333
$ws_ARM.Range("M38").Value = -2439.3333

# Row 102
$ws_ARM.Range("H102").Value = 501000
$ws_ARM.Range("I102").Value = 2000
$ws_ARM.Range("K102").Value = 2000
$ws_ARM.Range("M102").Value = -378

# Row 116
$ws_ARM.Range("H116").Value = 38656.63
$ws_ARM.Range("I116").Value = 1138.7
$ws_ARM.Range("J116").Value = 60726
$ws_ARM.Range("K116").Value = 1138.7
$ws_ARM.Range("L116").Value = 60726
$ws_ARM.Range("M116").Value = 1155.3
$ws_ARM.Range("N116").Value = -65314

# Row 132
$ws_ARM.Range("H132").Value = 6680.8
$ws_ARM.Range("I132").Value = 10715.23
$ws_ARM.Range("J132").Value = 3595.647
$ws_ARM.Range("K132").Value = 32145.69
$ws_ARM.Range("L132").Value = 10786.941
$ws_ARM.Range("M132").Value = -29615.69
$ws_ARM.Range("N132").Value = -15846.941

# --- BSM sheet ---
$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 3
$ws_BSM.Range("H3").Value = 38656.63
$ws_BSM.Range("I3").Value = 1138.7
$ws_BSM.Range("J3").Value = 60726
$ws_BSM.Range("K3").Value = 1138.7
$ws_BSM.Range("L3").Value = 60726
$ws_BSM.Range("M3").Value = -1024.7
$ws_BSM.Range("N3").Value = -60954

# Row 94
$ws_BSM.Range("H94").Value = 46605.547
$ws_BSM.Range("I94").Value = 518.6923
$ws_BSM.Range("J94").Value = 113175.445
$ws_BSM.Range("K94").Value = 518.6923
$ws_BSM.Range("L94").Value = 113175.445
$ws_BSM.Range("M94").Value = -67.69230000000005
$ws_BSM.Range("N94").Value = -114077.445

# Row 134
$ws_BSM.Range("H134").Value = 2236.4888
$ws_BSM.Range("I134").Value = 1914
$ws_BSM.Range("J134").Value = 3728
$ws_BSM.Range("K134").Value = 5742
$ws_BSM.Range("L134").Value = 11184
$ws_BSM.Range("M134").Value = -3207
$ws_BSM.Range("N134").Value = -16254

# --- CRP sheet ---
$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 107
$ws_CRP.Range("H107").Value = 418.5
$ws_CRP.Range("I107").Value = 322.53845
$ws_CRP.Range("J107").Value = 834.3333
$ws_CRP.Range("K107").Value = 322.53845
$ws_CRP.Range("L107").Value = 834.3333
$ws_CRP.Range("M107").Value = 1597.46155
$ws_CRP.Range("N107").Value = -4674.3333

# --- CUL sheet ---
$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 80
$ws_CUL.Range("H80").Value = 10111.444
$ws_CUL.Range("I80").Value = 16000
$ws_CUL.Range("K80").Value = 48000
$ws_CUL.Range("M80").Value = -47064

# Row 83
$ws_CUL.Range("H83").Value = 10111.444
$ws_CUL.Range("I83").Value = 16000
$ws_CUL.Range("K83").Value = 144000
$ws_CUL.Range("M83").Value = -139320

# --- GSM sheet ---
$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 40
$ws_GSM.Range("H40").Value = 50000
$ws_GSM.Range("J40").Value = 0
$ws_GSM.Range("L40").Value = 0
$ws_GSM.Range("N40").ClearContents()

# Row 107
$ws_GSM.Range("H107").Value = 1013.5263
$ws_GSM.Range("I107").Value = 1223.091
$ws_GSM.Range("J107").Value = 725.375
$ws_GSM.Range("K107").Value = 1223.091
$ws_GSM.Range("L107").Value = 725.375
$ws_GSM.Range("M107").Value = 696.9090000000001
$ws_GSM.Range("N107").Value = -4565.375

# Row 122
$ws_GSM.Range("H122").Value = 2661.7778
$ws_GSM.Range("I122").Value = 2244.8572
$ws_GSM.Range("K122").Value = 6734.571599999999
$ws_GSM.Range("M122").Value = -4284.571599999999

# Row 123
$ws_GSM.Range("H123").Value = 8622.799999999999
$ws_GSM.Range("J123").Value = 8622.799999999999
$ws_GSM.Range("L123").Value = 8622.799999999999
$ws_GSM.Range("N123").Value = -13522.8

# Row 126
$ws_GSM.Range("H126").Value = 2347.125
$ws_GSM.Range("I126").Value = 1667.5
$ws_GSM.Range("K126").Value = 5002.5
$ws_GSM.Range("M126").Value = -2532.5

# Row 132
$ws_GSM.Range("H132").Value = 2554.6956
$ws_GSM.Range("I132").Value = 2236
$ws_GSM.Range("J132").Value = 4330.2856
$ws_GSM.Range("K132").Value = 6708
$ws_GSM.Range("L132").Value = 12990.8568
$ws_GSM.Range("M132").Value = -4178
$ws_GSM.Range("N132").Value = -18050.8568

# --- LTW sheet ---
$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 136
$ws_LTW.Range("H136").Value = 18037422
$ws_LTW.Range("I136").Value = 23810810
$ws_LTW.Range("K136").Value = 71432430
$ws_LTW.Range("M136").Value = -71429880

# --- WVR sheet ---
$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 16
$ws_WVR.Range("H16").Value = 43472.668
$ws_WVR.Range("J16").Value = 43472.668
$ws_WVR.Range("L16").Value = 43472.668
$ws_WVR.Range("N16").Value = -44056.668

# Row 61
$ws_WVR.Range("H61").Value = 5367.3335
$ws_WVR.Range("I61").Value = 5367.3335
$ws_WVR.Range("K61").Value = 5367.3335
$ws_WVR.Range("M61").Value = -5075.3335

# Row 81
$ws_WVR.Range("H81").Value = 67005.8
$ws_WVR.Range("I81").Value = 68745.07000000001
$ws_WVR.Range("J81").Value = 61788
$ws_WVR.Range("K81").Value = 137490.14
$ws_WVR.Range("L81").Value = 123576
$ws_WVR.Range("M81").Value = -136429.14
$ws_WVR.Range("N81").Value = -125698

# Row 84
$ws_WVR.Range("H84").Value = 67005.8
$ws_WVR.Range("I84").Value = 68745.07000000001
$ws_WVR.Range("J84").Value = 61788
$ws_WVR.Range("K84").Value = 687450.7000000001
$ws_WVR.Range("L84").Value = 617880
$ws_WVR.Range("M84").Value = -682146.7000000001
$ws_WVR.Range("N84").Value = -628488

# Row 107
$ws_WVR.Range("H107").Value = 463.625
$ws_WVR.Range("I107").Value = 392.5
$ws_WVR.Range("J107").Value = 582.1667
$ws_WVR.Range("K107").Value = 1177.5
$ws_WVR.Range("L107").Value = 1746.5001
$ws_WVR.Range("M107").Value = 742.5
$ws_WVR.Range("N107").Value = -5586.5001
